$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted for "Cebollín" (row 319), which shifts every
# subsequent record down by one row (old row 319 becomes 320, ..., old row 454
# becomes the new row 455), growing the used range from A1:R454 to A1:R455.
$ws.Rows.Item(319).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Cells.Item(319, 1).Value  = 3
$ws.Cells.Item(319, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(319, 3).Value  = "Coquimbo"
$ws.Cells.Item(319, 4).Value  = 44636
$ws.Cells.Item(319, 5).Value  = 5
$ws.Cells.Item(319, 6).Value  = 100112037
$ws.Cells.Item(319, 7).Value  = "Cebollín"
$ws.Cells.Item(319, 8).Value  = "Sin especificar"
$ws.Cells.Item(319, 9).Value  = "Primera"
$ws.Cells.Item(319, 10).Value = 210
$ws.Cells.Item(319, 11).Value = 4000
$ws.Cells.Item(319, 12).Value = 4300
$ws.Cells.Item(319, 13).Value = 4129
$ws.Cells.Item(319, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(319, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(319, 16).Value = 115
$ws.Cells.Item(319, 17).Value = 36
$ws.Cells.Item(319, 18).Value = "Hortaliza"
